$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Bulk-update the "Förändrad" (changed) date column for every existing
#    data row (2..345) from 45190 (2023-09-21) to 45192 (2023-09-23).
$ws.Range("C2:C345").Value = 45192

# 2) Row 345 picks up an explicit row-height (matches every other data row).
$ws.Rows.Item(345).RowHeight = 15

# 3) Append the new notification as row 346.
$ws.Range("A346").Value = "A 44976-2023"

$ws.Range("B346").Value = 45190
$ws.Range("B346").NumberFormat = $ws.Range("B345").NumberFormat

$ws.Range("C346").Value = 45192
$ws.Range("C346").NumberFormat = $ws.Range("C345").NumberFormat

$ws.Range("D346").Value = "VÄSTERBOTTENS LÄN"
$ws.Range("E346").Value = "NORSJÖ"

$ws.Range("G346").Value = 2.1
$ws.Range("H346").Value = 0
$ws.Range("I346").Value = 0
$ws.Range("J346").Value = 0
$ws.Range("K346").Value = 0
$ws.Range("L346").Value = 0
$ws.Range("M346").Value = 0
$ws.Range("N346").Value = 0
$ws.Range("O346").Value = 0
$ws.Range("P346").Value = 0
$ws.Range("Q346").Value = 0

$ws.Range("R346").WrapText = $ws.Range("R345").WrapText
